$wb = $excel.ActiveWorkbook

# --- Rename sheet "burp" -> "Process" ---
$wsProcess = $wb.Worksheets.Item("burp")
$wsProcess.Name = "Process"

$wsInput = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# --- Input sheet: insert a new column G ("Visible") before the old G (Options) ---
$wsInput.Range("G1").EntireColumn.Insert()

$wsInput.Range("G1").Value = "Visible"
$wsInput.Range("G2").Value = $true
$wsInput.Range("G3").Value = $true
$wsInput.Range("G4").Value = $true
$wsInput.Range("G5").Value = $true

# --- Input sheet: fix bug, Width value D4 9 -> 8 ---
$wsInput.Range("D4").Value = 8

# --- Output sheet: update VLOOKUP formula to use renamed "Process" sheet ---
$wsOutput.Range("C5").Formula = "=VLOOKUP(Input!D5,Process!A2:B5,2,FALSE)*C3/1000"

# --- Selections / active sheet, matching the authored view state ---
$wsProcess.Range("B6").Select() | Out-Null
$wsOutput.Range("C6").Select() | Out-Null
$wsInput.Range("G6").Select() | Out-Null
$wsInput.Activate() | Out-Null

Write-Host "edit applied"
